$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H54").Value = 40084
$ws.Range("J54").Value = 40084
$ws.Range("L54").Value = 40084
$ws.Range("N54").Value = -41056
$ws.Range("H62").Value = 1825.85
$ws.Range("I62").Value = 1697.7693
$ws.Range("J62").Value = 2063.7144
$ws.Range("K62").Value = 1697.7693
$ws.Range("L62").Value = 2063.7144
$ws.Range("M62").Value = -1073.7693
$ws.Range("N62").Value = -3311.7144
$ws.Range("H65").Value = 1825.85
$ws.Range("I65").Value = 1697.7693
$ws.Range("J65").Value = 2063.7144
$ws.Range("K65").Value = 8488.8465
$ws.Range("L65").Value = 10318.572
$ws.Range("M65").Value = -5368.8465
$ws.Range("N65").Value = -16558.572
$ws.Range("H98").Value = 1207.7778
$ws.Range("I98").Value = 1025.6666
$ws.Range("K98").Value = 1025.6666
$ws.Range("M98").Value = 472.3334
$ws.Range("H100").Value = 1324.4
$ws.Range("I100").Value = 1083.2727
$ws.Range("K100").Value = 1083.2727
$ws.Range("M100").Value = -542.2727
$ws.Range("H122").Value = 1207.7778
$ws.Range("I122").Value = 1025.6666
$ws.Range("K122").Value = 3076.9998
$ws.Range("M122").Value = -626.9998000000001
$ws.Range("H125").Value = 1066.6666
$ws.Range("I125").Value = 600
$ws.Range("J125").Value = 2000
$ws.Range("K125").Value = 5400
$ws.Range("L125").Value = 18000
$ws.Range("M125").Value = -2940
$ws.Range("N125").Value = -22920
$ws.Range("H138").Value = 4304.8
$ws.Range("I138").Value = 6083.909
$ws.Range("J138").Value = 3489.375
$ws.Range("K138").Value = 18251.727
$ws.Range("L138").Value = 10468.125
$ws.Range("M138").Value = -13111.727
$ws.Range("N138").Value = -20748.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 6635.4287
$ws.Range("I46").Value = 3874.75
$ws.Range("K46").Value = 3874.75
$ws.Range("M46").Value = -3555.75
$ws.Range("H61").Value = 2313.1628
$ws.Range("I61").Value = 2174.7693
$ws.Range("J61").Value = 3662.5
$ws.Range("K61").Value = 2174.7693
$ws.Range("L61").Value = 3662.5
$ws.Range("M61").Value = -1962.7693
$ws.Range("N61").Value = -4086.5
$ws.Range("H110").Value = 1415.4445
$ws.Range("I110").Value = 1380.2941
$ws.Range("K110").Value = 1380.2941
$ws.Range("M110").Value = 664.7058999999999
$ws.Range("H122").Value = 1989.4445
$ws.Range("I122").Value = 2050.8333
$ws.Range("K122").Value = 6152.499899999999
$ws.Range("M122").Value = -3702.499899999999
$ws.Range("H136").Value = 2313.1628
$ws.Range("I136").Value = 2174.7693
$ws.Range("J136").Value = 3662.5
$ws.Range("K136").Value = 6524.3079
$ws.Range("L136").Value = 10987.5
$ws.Range("M136").Value = -3974.3079
$ws.Range("N136").Value = -16087.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 70492.5
$ws.Range("J13").Value = 70492.5
$ws.Range("L13").Value = 70492.5
$ws.Range("N13").Value = -70828.5
$ws.Range("H94").Value = 6606.346
$ws.Range("I94").Value = 6493.9443
$ws.Range("J94").Value = 6859.25
$ws.Range("K94").Value = 6493.9443
$ws.Range("L94").Value = 6859.25
$ws.Range("M94").Value = -6042.9443
$ws.Range("N94").Value = -7761.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2419.7058
$ws.Range("I31").Value = 3511.125
$ws.Range("K31").Value = 3511.125
$ws.Range("M31").Value = -3216.125
$ws.Range("H34").Value = 2419.7058
$ws.Range("I34").Value = 3511.125
$ws.Range("K34").Value = 3511.125
$ws.Range("M34").Value = -3309.125
$ws.Range("H99").Value = 11907.333
$ws.Range("I99").Value = 17543.77
$ws.Range("K99").Value = 17543.77
$ws.Range("M99").Value = -16045.77
$ws.Range("H122").Value = 11707.77
$ws.Range("I122").Value = 2278.0557
$ws.Range("J122").Value = 32924.625
$ws.Range("K122").Value = 6834.1671
$ws.Range("L122").Value = 98773.875
$ws.Range("M122").Value = -4384.1671
$ws.Range("N122").Value = -103673.875
$ws.Range("H126").Value = 11907.333
$ws.Range("I126").Value = 17543.77
$ws.Range("K126").Value = 52631.31
$ws.Range("M126").Value = -50161.31

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 59782.61
$ws.Range("J37").Value = 59782.61
$ws.Range("L37").Value = 179347.83
$ws.Range("N37").Value = -179571.83
$ws.Range("H44").Value = 6950
$ws.Range("J44").Value = 7611.1113
$ws.Range("L44").Value = 22833.3339
$ws.Range("N44").Value = -23629.3339
$ws.Range("H62").Value = 3550
$ws.Range("J62").Value = 4066.6667
$ws.Range("L62").Value = 12200.0001
$ws.Range("N62").Value = -13572.0001
$ws.Range("H65").Value = 3550
$ws.Range("J65").Value = 4066.6667
$ws.Range("L65").Value = 36600.0003
$ws.Range("N65").Value = -43464.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 70
$ws.Range("I9").Value = 70
$ws.Range("K9").Value = 70
$ws.Range("M9").Value = 100
$ws.Range("H122").Value = 3483.125
$ws.Range("I122").Value = 3517.6924
$ws.Range("K122").Value = 10553.0772
$ws.Range("M122").Value = -8103.0772
$ws.Range("H126").Value = 2471.25
$ws.Range("I126").Value = 2005.5834
$ws.Range("K126").Value = 6016.7502
$ws.Range("M126").Value = -3546.7502

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5105.727
$ws.Range("I7").Value = 4816.8
$ws.Range("J7").Value = 7995
$ws.Range("K7").Value = 4816.8
$ws.Range("L7").Value = 7995
$ws.Range("M7").Value = -4704.8
$ws.Range("N7").Value = -8219
$ws.Range("H9").Value = 5088.75
$ws.Range("I9").Value = 177.5
$ws.Range("K9").Value = 177.5
$ws.Range("M9").Value = 46.5
$ws.Range("H30").Value = 1020
$ws.Range("I30").Value = 1000
$ws.Range("J30").Value = 1100
$ws.Range("K30").Value = 1000
$ws.Range("L30").Value = 1100
$ws.Range("M30").Value = -892
$ws.Range("N30").Value = -1316
$ws.Range("H35").Value = 25258.5
$ws.Range("I35").Value = 5499.5
$ws.Range("J35").Value = 45017.5
$ws.Range("K35").Value = 5499.5
$ws.Range("L35").Value = 45017.5
$ws.Range("M35").Value = -5163.5
$ws.Range("N35").Value = -45689.5
$ws.Range("H61").Value = 7902.6343
$ws.Range("I61").Value = 7982
$ws.Range("J61").Value = 7517.143
$ws.Range("K61").Value = 7982
$ws.Range("L61").Value = 7517.143
$ws.Range("M61").Value = -7780
$ws.Range("N61").Value = -7921.143
$ws.Range("H99").Value = 44947
$ws.Range("I99").Value = 44947
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 44947
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -41952
$ws.Range("H113").Value = 7902.6343
$ws.Range("I113").Value = 7982
$ws.Range("J113").Value = 7517.143
$ws.Range("K113").Value = 7982
$ws.Range("L113").Value = 7517.143
$ws.Range("M113").Value = -5812
$ws.Range("N113").Value = -11857.143
$ws.Range("H122").Value = 5070.8
$ws.Range("I122").Value = 4339.25
$ws.Range("K122").Value = 13017.75
$ws.Range("M122").Value = -10567.75
$ws.Range("H126").Value = 5105.727
$ws.Range("I126").Value = 4816.8
$ws.Range("J126").Value = 7995
$ws.Range("K126").Value = 14450.4
$ws.Range("L126").Value = 23985
$ws.Range("M126").Value = -11980.4
$ws.Range("N126").Value = -28925
$ws.Range("H132").Value = 3773.0588
$ws.Range("I132").Value = 3242.4614
$ws.Range("K132").Value = 9727.3842
$ws.Range("M132").Value = -7197.3842
$ws.Range("H136").Value = 2645.3784
$ws.Range("I136").Value = 2099.7334
$ws.Range("J136").Value = 4983.857
$ws.Range("K136").Value = 6299.2002
$ws.Range("L136").Value = 14951.571
$ws.Range("M136").Value = -3749.2002
$ws.Range("N136").Value = -20051.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 26820.924
$ws.Range("I41").Value = 7499
$ws.Range("J41").Value = 28431.084
$ws.Range("K41").Value = 7499
$ws.Range("L41").Value = 28431.084
$ws.Range("M41").Value = -7109
$ws.Range("N41").Value = -29211.084
$ws.Range("H122").Value = 42113.9
$ws.Range("I122").Value = 3499.4348
$ws.Range("J122").Value = 190136
$ws.Range("K122").Value = 10498.3044
$ws.Range("L122").Value = 570408
$ws.Range("M122").Value = -8048.304400000001
$ws.Range("N122").Value = -575308
$ws.Range("H126").Value = 1489.037
$ws.Range("I126").Value = 1307.9546
$ws.Range("K126").Value = 3923.8638
$ws.Range("M126").Value = -1453.8638
$ws.Range("H136").Value = 1305.7407
$ws.Range("I136").Value = 1138.8334
$ws.Range("J136").Value = 1889.9166
$ws.Range("K136").Value = 3416.5002
$ws.Range("L136").Value = 5669.7498
$ws.Range("M136").Value = -866.5001999999999
$ws.Range("N136").Value = -10769.7498
